$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 832, shifting existing rows 832:877 down to 833:878.
$ws.Rows.Item(832).Insert()

# Populate the newly inserted row 832 with the new weekly record.
$ws.Cells.Item(832, 1).Value = 4
$ws.Cells.Item(832, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(832, 3).Value = "Los Lagos"
$ws.Cells.Item(832, 4).Value = 45267
$ws.Cells.Item(832, 5).Value = 10
$ws.Cells.Item(832, 6).Value = 100112006
$ws.Cells.Item(832, 7).Value = "Repollo"
$ws.Cells.Item(832, 8).Value = "Copenhague"
$ws.Cells.Item(832, 9).Value = "Primera"
$ws.Cells.Item(832, 10).Value = 600
$ws.Cells.Item(832, 11).Value = 1900
$ws.Cells.Item(832, 12).Value = 2000
$ws.Cells.Item(832, 13).Value = 1950
$ws.Cells.Item(832, 14).Value = "`$/unidad"
$ws.Cells.Item(832, 15).Value = "Región Metropolitana"
$ws.Cells.Item(832, 16).Value = 1950
$ws.Cells.Item(832, 17).Value = 1
$ws.Cells.Item(832, 18).Value = "Hortaliza"
